# Generate Report for Handoff
#
# A new handoff run produced a fresh commit GUID and content hash, and the
# localization-status report needs to reflect it: the Overview sheet's
# source-file link/name and "Latest Handoff Date", plus each per-locale
# sheet's (zh-cn, de-de) source link, generated .xlf link, and handoff
# timestamp.

$wb = $excel.ActiveWorkbook

$oldGuid = "2bc5cede-3bd4-4d25-9f6f-6b70eb0fcbea"
$newGuid = "3f7af0e6-3842-406e-91a3-430a4d9c9fb8"
$oldHash = "c936227c04c747f2440d2b0b013210d37197cc37"
$newHash = "a63d2eca2096107008b7e2821824aafe699d34d0"

$oldMdName = "$oldGuid.md"
$newMdName = "$newGuid.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = "2016-50-11 10:50:21"

foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMdName) {
        $hl.TextToDisplay = $newMdName
    }
}

# --- zh-cn sheet ---
$oldZhXlf = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = "2016-03-11 10:50:19"

foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMdName) {
        $hl.TextToDisplay = $newMdName
    } elseif ($hl.TextToDisplay -eq $oldZhXlf) {
        $hl.TextToDisplay = $newZhXlf
    }
}

# --- de-de sheet ---
$oldDeXlf = "$oldGuid.$oldHash.de-de.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = "2016-03-11 10:50:21"

foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMdName) {
        $hl.TextToDisplay = $newMdName
    } elseif ($hl.TextToDisplay -eq $oldDeXlf) {
        $hl.TextToDisplay = $newDeXlf
    }
}
